$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# "textures\grnglow.igb" (row 51) is an unused texture that doesn't exist
# anywhere and isn't the default for any version of the game - remove it.
$ws.Rows.Item(51).Delete()

# The conditional-formatting rules covered the whole table (A2:A123,
# B2:B123, ... F2:F123). After the row deletion the table is one row
# shorter, so re-point each rule's range at the new extent (…122).
$lastRow = 122

$fcs = $ws.Range("A2:A" + $lastRow).FormatConditions
for ($i = 1; $i -le $fcs.Count; $i++) {
    $fcs.Item($i).ModifyAppliesToRange($ws.Range("A2:A" + $lastRow))
}

$fcs = $ws.Range("B2:B" + $lastRow).FormatConditions
for ($i = 1; $i -le $fcs.Count; $i++) {
    $fcs.Item($i).ModifyAppliesToRange($ws.Range("B2:B" + $lastRow))
}

$fcs = $ws.Range("B2:F" + $lastRow).FormatConditions
for ($i = 1; $i -le $fcs.Count; $i++) {
    $fcs.Item($i).ModifyAppliesToRange($ws.Range("B2:F" + $lastRow))
}

$fcs = $ws.Range("C2:C" + $lastRow).FormatConditions
for ($i = 1; $i -le $fcs.Count; $i++) {
    $fcs.Item($i).ModifyAppliesToRange($ws.Range("C2:C" + $lastRow))
}

$fcs = $ws.Range("D2:D" + $lastRow).FormatConditions
for ($i = 1; $i -le $fcs.Count; $i++) {
    $fcs.Item($i).ModifyAppliesToRange($ws.Range("D2:D" + $lastRow))
}

$fcs = $ws.Range("E2:E" + $lastRow).FormatConditions
for ($i = 1; $i -le $fcs.Count; $i++) {
    $fcs.Item($i).ModifyAppliesToRange($ws.Range("E2:E" + $lastRow))
}

$fcs = $ws.Range("F2:F" + $lastRow).FormatConditions
for ($i = 1; $i -le $fcs.Count; $i++) {
    $fcs.Item($i).ModifyAppliesToRange($ws.Range("F2:F" + $lastRow))
}

# Keep the selection/scroll position close to where the row was removed.
$ws.Application.ActiveWindow.ScrollRow = 31
$ws.Range("E39").Select()
